# Boolean Do Suppliers Bid at Peak Capacity Factors - update to v2.0.0
# (adds crude oil / heavy-or-residual-fuel-oil / municipal-solid-waste rows,
#  flips a few booleans, and relabels/restyles the header cell)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BDSBaPCF")

# --- Flip a few existing boolean flags ---
$ws.Range("B4").Value = 1    # nuclear
$ws.Range("B11").Value = 0   # petroleum
$ws.Range("B12").Value = 0   # natural gas peaker

# --- New fuel-type rows, mirroring existing boolean values via formula ---
$ws.Range("A15").Value = "crude oil"
$ws.Range("B15").Formula = "=B11"

$ws.Range("A16").Value = "heavy or residual fuel oil"
$ws.Range("B16").Formula = "=B11"

$ws.Range("A17").Value = "municipal solid waste"
$ws.Range("B17").Formula = "=B9"

# --- Header cell B1: new label text + bold/wrap/right-aligned style ---
$header = $ws.Range("B1")
$header.Value = "Do Suppliers Bid at Peak Capacity Factors (Boolean)"
$header.Font.Bold = $true
$header.WrapText = $true
$header.HorizontalAlignment = -4152
$ws.Rows.Item(1).RowHeight = 45

$ws.Range("B1").Select()
